$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.249.24'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.650.86'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.514'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.37%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.257'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.40%  '
$ws.Range('E9').Value = '  +0.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.07'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.68%  '
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.882.73'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.662.21'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.14'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.542'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.65'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.241.27'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '220.54'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.88%  '
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.56'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.27%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.79'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('B26').Value = 'BinanceUSD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.53'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.84'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0509'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.03%  '
$ws.Range('E31').Value = '  -0.69%  '
$ws.Range('E32').Value = '  -0.96%  '
$ws.Range('E33').Value = '  +0.59%  '
$ws.Range('E34').Value = '  +0.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.259.97'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E37').Value = '  -0.26%  '
$ws.Range('E38').Value = '  +1.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.844'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.47%  '
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.44'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.21'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.792.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.52%  '
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.73'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.60'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0106'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +23.48%  '
$ws.Range('E49').Value = '  -0.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.68'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0970'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.26%  '
